# Updates computed-profit columns (H:N) across the leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to refreshed market-price
# derived figures. A couple of rows (ALC!N98, ALC!N122) had their
# trailing column dropped entirely, so those are cleared rather than
# rewritten.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 696.875
$ws.Range("I4").Value = 344.5
$ws.Range("J4").Value = 1049.25
$ws.Range("K4").Value = 344.5
$ws.Range("L4").Value = 1049.25
$ws.Range("M4").Value = -230.5
$ws.Range("N4").Value = -1277.25
$ws.Range("H17").Value = 752.2941
$ws.Range("J17").Value = 752.2941
$ws.Range("L17").Value = 2256.8823
$ws.Range("N17").Value = -2592.8823
$ws.Range("H40").Value = 25226.908
$ws.Range("I40").Value = 7416.3335
$ws.Range("K40").Value = 7416.3335
$ws.Range("M40").Value = -7241.3335
$ws.Range("H43").Value = 14665
$ws.Range("I43").Value = 18243.5
$ws.Range("J43").Value = 10370.8
$ws.Range("K43").Value = 18243.5
$ws.Range("L43").Value = 10370.8
$ws.Range("M43").Value = -18174.5
$ws.Range("N43").Value = -10508.8
$ws.Range("H86").Value = 4316.125
$ws.Range("I86").Value = 4234.5
$ws.Range("K86").Value = 4234.5
$ws.Range("M86").Value = -3111.5
$ws.Range("H88").Value = 5170.1816
$ws.Range("I88").Value = 5148.75
$ws.Range("J88").Value = 5182.4287
$ws.Range("K88").Value = 5148.75
$ws.Range("L88").Value = 5182.4287
$ws.Range("M88").Value = -4742.75
$ws.Range("N88").Value = -5994.4287
$ws.Range("H89").Value = 4316.125
$ws.Range("I89").Value = 4234.5
$ws.Range("K89").Value = 21172.5
$ws.Range("M89").Value = -15556.5
$ws.Range("H91").Value = 5170.1816
$ws.Range("I91").Value = 5148.75
$ws.Range("J91").Value = 5182.4287
$ws.Range("K91").Value = 5148.75
$ws.Range("L91").Value = 5182.4287
$ws.Range("M91").Value = -3744.75
$ws.Range("N91").Value = -7990.4287
$ws.Range("H92").Value = 820.86664
$ws.Range("I92").Value = 641.8461
$ws.Range("J92").Value = 1984.5
$ws.Range("K92").Value = 641.8461
$ws.Range("L92").Value = 1984.5
$ws.Range("M92").Value = 606.1539
$ws.Range("N92").Value = -4480.5
$ws.Range("M98").Value = -576.25
$ws.Range("H98").Value = 2074.25
$ws.Range("I98").Value = 2074.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2074.25
$ws.Range("L98").Value = 0
$ws.Range("H113").Value = 4950.222
$ws.Range("I113").Value = 4744.25
$ws.Range("J113").Value = 5115
$ws.Range("K113").Value = 4744.25
$ws.Range("L113").Value = 5115
$ws.Range("M113").Value = -1490.25
$ws.Range("N113").Value = -11623
$ws.Range("M122").Value = -3772.75
$ws.Range("H122").Value = 2074.25
$ws.Range("I122").Value = 2074.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6222.75
$ws.Range("L122").Value = 0
$ws.Range("H135").Value = 7376.222
$ws.Range("I135").Value = 1415
$ws.Range("K135").Value = 12735
$ws.Range("M135").Value = -10200
$ws.Range("H137").Value = 2083.889
$ws.Range("I137").Value = 1056.4
$ws.Range("K137").Value = 3169.2
$ws.Range("M137").Value = -619.2000000000003
$ws.Range("H138").Value = 3580.394
$ws.Range("I138").Value = 1797.2
$ws.Range("J138").Value = 4355.696
$ws.Range("K138").Value = 5391.6
$ws.Range("L138").Value = 13067.088
$ws.Range("M138").Value = -251.6000000000004
$ws.Range("N138").Value = -23347.088
$ws.Range("N98").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1339.8591
$ws.Range("I32").Value = 781.3677
$ws.Range("K32").Value = 781.3677
$ws.Range("M32").Value = -494.3677
$ws.Range("H97").Value = 1367.75
$ws.Range("I97").Value = 1367.75
$ws.Range("K97").Value = 1367.75
$ws.Range("M97").Value = -871.75
$ws.Range("H102").Value = 7829.5
$ws.Range("I102").Value = 7435.6665
$ws.Range("K102").Value = 7435.6665
$ws.Range("M102").Value = -5813.6665
$ws.Range("H110").Value = 2829.077
$ws.Range("I110").Value = 2099.375
$ws.Range("K110").Value = 2099.375
$ws.Range("M110").Value = -54.375
$ws.Range("H132").Value = 3305.9546
$ws.Range("I132").Value = 2986
$ws.Range("K132").Value = 8958
$ws.Range("M132").Value = -6428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2925.4285
$ws.Range("I80").Value = 2872.25
$ws.Range("J80").Value = 2996.3333
$ws.Range("K80").Value = 2872.25
$ws.Range("L80").Value = 2996.3333
$ws.Range("M80").Value = -1874.25
$ws.Range("N80").Value = -4992.3333
$ws.Range("H83").Value = 2925.4285
$ws.Range("I83").Value = 2872.25
$ws.Range("J83").Value = 2996.3333
$ws.Range("K83").Value = 14361.25
$ws.Range("L83").Value = 14981.6665
$ws.Range("M83").Value = -9369.25
$ws.Range("N83").Value = -24965.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2102.8667
$ws.Range("I105").Value = 2475.7
$ws.Range("K105").Value = 2475.7
$ws.Range("M105").Value = -728.6999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 18287.143
$ws.Range("I70").Value = 12005
$ws.Range("K70").Value = 36015
$ws.Range("M70").Value = -35700
$ws.Range("H73").Value = 18287.143
$ws.Range("I73").Value = 12005
$ws.Range("K73").Value = 36015
$ws.Range("M73").Value = -34923
$ws.Range("H75").Value = 2699.6667
$ws.Range("I75").Value = 2461.25
$ws.Range("J75").Value = 2818.875
$ws.Range("K75").Value = 7383.75
$ws.Range("L75").Value = 8456.625
$ws.Range("M75").Value = -6385.75
$ws.Range("N75").Value = -10452.625
$ws.Range("H78").Value = 2699.6667
$ws.Range("I78").Value = 2461.25
$ws.Range("J78").Value = 2818.875
$ws.Range("K78").Value = 22151.25
$ws.Range("L78").Value = 25369.875
$ws.Range("M78").Value = -17159.25
$ws.Range("N78").Value = -35353.875
$ws.Range("H94").Value = 15859
$ws.Range("J94").Value = 15835.429
$ws.Range("L94").Value = 47506.287
$ws.Range("N94").Value = -48858.287
$ws.Range("H100").Value = 13332.333
$ws.Range("J100").Value = 13332.333
$ws.Range("L100").Value = 39996.999
$ws.Range("N100").Value = -41618.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 819.6667
$ws.Range("I97").Value = 826.3333
$ws.Range("J97").Value = 813
$ws.Range("K97").Value = 826.3333
$ws.Range("L97").Value = 813
$ws.Range("M97").Value = -330.3333
$ws.Range("N97").Value = -1805
$ws.Range("H111").Value = 19999.5
$ws.Range("J111").Value = 19999.5
$ws.Range("L111").Value = 19999.5
$ws.Range("N111").Value = -26133.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15002.223
$ws.Range("I7").Value = 18183.166
$ws.Range("J7").Value = 8640.333000000001
$ws.Range("K7").Value = 18183.166
$ws.Range("L7").Value = 8640.333000000001
$ws.Range("M7").Value = -18071.166
$ws.Range("N7").Value = -8864.333000000001
$ws.Range("H16").Value = 831.7895
$ws.Range("I16").Value = 683.73334
$ws.Range("J16").Value = 1387
$ws.Range("K16").Value = 683.73334
$ws.Range("L16").Value = 1387
$ws.Range("M16").Value = -513.73334
$ws.Range("N16").Value = -1727
$ws.Range("H40").Value = 14012.571
$ws.Range("I40").Value = 16813.143
$ws.Range("J40").Value = 11212
$ws.Range("K40").Value = 16813.143
$ws.Range("L40").Value = 11212
$ws.Range("M40").Value = -16677.143
$ws.Range("N40").Value = -11484
$ws.Range("H100").Value = 4539.3076
$ws.Range("I100").Value = 1703.6666
$ws.Range("K100").Value = 1703.6666
$ws.Range("M100").Value = -1162.6666
$ws.Range("H126").Value = 15002.223
$ws.Range("I126").Value = 18183.166
$ws.Range("J126").Value = 8640.333000000001
$ws.Range("K126").Value = 54549.49800000001
$ws.Range("L126").Value = 25920.999
$ws.Range("M126").Value = -52079.49800000001
$ws.Range("N126").Value = -30860.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 73740
$ws.Range("J95").Value = 73740
$ws.Range("L95").Value = 73740
$ws.Range("N95").Value = -79232
$ws.Range("H96").Value = 2446.3333
$ws.Range("J96").Value = 2435.6
$ws.Range("L96").Value = 2435.6
$ws.Range("N96").Value = -5181.6
$ws.Range("H100").Value = 2220.8333
$ws.Range("J100").Value = 2265
$ws.Range("L100").Value = 4530
$ws.Range("N100").Value = -5612
$ws.Range("H136").Value = 2290.9211
$ws.Range("I136").Value = 1740.3793
$ws.Range("J136").Value = 4064.889
$ws.Range("K136").Value = 5221.1379
$ws.Range("L136").Value = 12194.667
$ws.Range("M136").Value = -2671.1379
$ws.Range("N136").Value = -17294.667
